$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix typo in header B1 ---
$ws.Range("B1").Value = "spatial resolution"

# --- Add new header for column H ---
$ws.Range("H1").Value = "comment"

# --- Row 2 (ERA5): change spatial resolution from numeric 0.5 to text "0.5 degree",
#     and temporal resolution from "daily" to "daily (mean)" ---
$ws.Range("B2").Value = "0.5 degree"
$ws.Range("C2").Value = "daily (mean)"
# D2, E2, F2, G2 (time period, external source, local location, file format) stay the same values

# --- Row 3 (Barentswatch) - new row ---
$ws.Range("A3").Value = "Barentswatch"
$ws.Range("B3").Value = "point measurements"
$ws.Range("C3").Value = "weekly (mean)"
$ws.Range("D3").Value = "2012-2020(2021?)"
$ws.Range("E3").Value = "Provided my the fishfarms"
$ws.Range("F3").Value = "no location"
$ws.Range("G3").Value = "csv"
$ws.Range("H3").Value = "unsupervised post processing done by provider, missing data/timelag might be an issue."

# --- Row 4 (RTG_SST_HR) - new row ---
$ws.Range("A4").Value = "RTG_SST_HR"
$ws.Range("B4").Value = "1/12 degree"
$ws.Range("C4").Value = "daily"
$ws.Range("D4").Value = "2005/09-2020/02"
$ws.Range("E4").Value = "NOAA"
$ws.Range("F4").Value = "no location"
$ws.Range("G4").Value = "netCDF"
$ws.Range("H4").Value = "https://polar.ncep.noaa.gov/sst/ophi/ no longer available?"

# Row 4 has a slightly reduced height vs the default
$ws.Rows.Item(4).RowHeight = 15.65

# --- Hyperlink on H4 pointing to the NOAA page ---
$ws.Hyperlinks.Add($ws.Range("H4"), "https://polar.ncep.noaa.gov/sst/ophi/")

# --- Update selection to H4 (matches author's last-edited cell) ---
[void]$ws.Range("H4").Select()

Write-Output "done"
